$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-8 from 45224 to 45233
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45233
}
